$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Week 2")
$ws.Select()

$ws.Range("C10").Select()

$ws.Range("C8").Value = 0.57638888888888895
$ws.Range("D8").Value = 0.60069444444444442

$ws.Range("C9").Value = 0.63541666666666663
$ws.Range("D9").Value = 0.73611111111111116
$ws.Range("E9").Value = 15

$ws.Range("G7:G10").Merge()
$ws.Range("H7:H10").Merge()

$ws.Range("G7:G10").HorizontalAlignment = -4108
$ws.Range("G7:G10").VerticalAlignment = -4108
$ws.Range("H7:H10").HorizontalAlignment = -4108
$ws.Range("H7:H10").VerticalAlignment = -4108
$ws.Range("H7:H10").WrapText = $true
